$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# ExcludeFields column repurposed as "Csvson" (CSV->JSON low-code column)
$ws.Range("J1").Value = "Csvson"
# New trailing column header for variable substitution support
$ws.Range("N1").Value = "AddifyVariables"

# --- Existing PetPost row (row 2) ---
# Clear the (now unused / relocated) ResponseByFields value, keep formatting
$ws.Range("G2").ClearContents()

# --- New PetGet row (row 3) demonstrating the low-code example ---
$ws.Range("A3").Value = "PetGet"
$ws.Range("B3").Value = "get by Id"
$ws.Range("C3").Value = "https://live.virtualandemo.com/api/pets/[petId]`n"
$ws.Range("D3").Value = "application/json"
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 12
$ws.Range("G3").Value = "name=Rocky`n"
$ws.Range("I3").Value = "GET"
$ws.Range("J3").Value = "id,name, category/id:name,status`ni~101,Rockey,i~100:german shepherd,available`n"
$ws.Range("K3").Value = 200
$ws.Range("N3").Value = "petId=100"

# Cells above contain embedded line breaks; re-fit the row so it keeps the
# sheet's default (non-custom) height rather than auto-growing for wrap.
$ws.Rows.Item(3).AutoFit()

# --- Column C width: author narrowed it considerably ---
$ws.Columns.Item(3).ColumnWidth = 42.8

# --- Selection moved to K3 ---
$ws.Range("K3").Select() | Out-Null
